# Including results from some generic routers
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of device data: HOSTNAME, STIG, VENDOR, TYPE
$data = @(
    @("r1",  "U_Network_Perimeter_Router_Cisco_STIG_V8R26_Manual-xccdf.xml.stig",       "CISCO", "PERIMETER"),
    @("r2",  "U_Network_Infrastructure_Router_Cisco_STIG_V8R23_Manual-xccdf.xml.stig",  "CISCO", "ROUTER"),
    @("r3",  "U_Network_Infrastructure_Router_Cisco_STIG_V8R23_Manual-xccdf.xml.stig",  "CISCO", "ROUTER"),
    @("sw1", "U_Network_L2_Switch_Cisco_STIG_V8R21_Manual-xccdf.xml.stig",              "CISCO", "L2_SWITCH"),
    @("sw2", "U_Network_L2_Switch_Cisco_STIG_V8R21_Manual-xccdf.xml.stig",              "CISCO", "L2_SWITCH")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $row++
}
